$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '24.772.54'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '1.701.62'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '''317.02'
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = '''0.3949'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = '''0.4083'
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("D9").Value = '''1.513'
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D11").Value = '''53.27'
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").Value = '''0.08915'
$ws.Range("D13").Value = '''7.722'
$ws.Range("E13").Value = '  +6.94%  '
$ws.Range("D14").Value = '''23.81'
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = '''8.232'
$ws.Range("E15").Value = '  +5.73%  '
$ws.Range("D16").Value = '''0.00001329'
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").Value = '1.702.19'
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '''99.69'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").Value = '''0.07135'
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").Value = '''19.97'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = '''7.151'
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("D22").Value = '''1.006'
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("D23").Value = '''14.81'
$ws.Range("E23").Value = '  +5.35%  '
$ws.Range("D24").Value = '24.759.66'
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").Value = '''3.152'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").Value = '''2.349'
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '''23.02'
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("D28").Value = '''9.280'
$ws.Range("E28").Value = '  +21.81%  '
$ws.Range("D29").Value = '''164.69'
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").Value = '''138.90'
$ws.Range("E30").Value = '  +3.18%  '
$ws.Range("D31").Value = '''5.158'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").Value = '''7.943'
$ws.Range("E32").Value = '  +7.88%  '
$ws.Range("D33").Value = '''0.09063'
$ws.Range("E33").Value = '  +6.06%  '
$ws.Range("D34").Value = '''1.074'
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("D35").Value = '''0.03043'
$ws.Range("E35").Value = '  +10.86%  '
$ws.Range("D36").Value = '''0.2792'
$ws.Range("E37").Value = '  -2.04%  '
$ws.Range("D38").Value = '''1.957'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '''14.49'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").Value = '''0.09302'
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''1.477'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.7805'
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("D43").Value = '''16.11'
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("D44").Value = '''2.647'
$ws.Range("E44").Value = '  +4.94%  '
$ws.Range("D45").Value = '''0.7255'
$ws.Range("E45").Value = '  +1.18%  '
$ws.Range("D46").Value = '''4.239'
$ws.Range("E46").Value = '  +0.75%  '
$ws.Range("D47").Value = '''1.357'
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '''140.48'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").Value = '''0.08011'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").Value = '''90.50'
$ws.Range("E51").Value = '  +2.71%  '